$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "Step 2" text rewritten + its expected result
$ws.Range("C3").Value = 'Step 2: Log in as a user with the appropriate role '
$ws.Range("D3").Value = 'I am redirected to the user''s dashboard'

# Row 4: "Step 3" text rewritten + its expected result
$ws.Range("C4").Value = 'Step 3: Go to the "Team Assessments" page '
$ws.Range("D4").Value = 'I am redirected to a page that has a list of assessments of people under my team '

# Row 5: "Step 4" text rewritten + its expected result
$ws.Range("C5").Value = 'Step 4: Try to delete an assessment of someone under my team'
$ws.Range("D5").Value = 'The assessment is removed from the database'

# Row 6 (new content): "Step 5" + expected result
$ws.Range("C6").Value = 'Step 5: Try to delete an assessment about the current logged in user'
$ws.Range("D6").Value = 'I am denied access to this'

# Row 7 (new content): "Step 6" (renumbered) + expected result
$ws.Range("C7").Value = 'Step 6: Try to delete an assessment of someone who is not under my team'
$ws.Range("D7").Value = 'I am denied access to this'

# Scroll the view down a bit and select C6, matching the author's final cursor position
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("C6").Select()
